# Adds a new "Git commit" slide (sentencia 4 git commit y creacion de rama 4)
# at the end of the deck, using the same "Título y objetos" (Title and
# Content) layout as the existing slides 2-4.

$p = $ppt.ActivePresentation

# ppLayoutText (2) == the "Título y objetos" custom layout already used by
# the other content slides (slide2.xml - slide4.xml in this deck).
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# --- Title placeholder: "Git commit" (centered) ---
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Git"
$title.InsertAfter(" ") | Out-Null
$title.InsertAfter("commit") | Out-Null
$title.ParagraphFormat.Alignment = 2  # ppAlignCenter

# --- Content placeholder: explanatory paragraph ---
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "El comando "
$body.InsertAfter("git") | Out-Null
$body.InsertAfter(" ") | Out-Null
$body.InsertAfter("commit") | Out-Null
$body.InsertAfter(" toma todos los contenidos de los archivos a los que se les realiza el seguimiento con ") | Out-Null
$body.InsertAfter("git") | Out-Null
$body.InsertAfter(" ") | Out-Null
$body.InsertAfter("add") | Out-Null
$body.InsertAfter(" y registra una nueva instantánea permanente en la base de datos y luego avanza el puntero de la rama en la rama actual.") | Out-Null
